# Atualizado por script em 11-11-2023 14:45
# Appends the new match row (index 76 / sheet row 77) for
# Pyunik Yerevan vs Noah (armenia / premier-league / 2023-2024).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 77
$prevRow = $newRow - 1

# Copy formatting from the row above first (keeps the existing cellXf
# for the index column (bold/border/center) and the date-time column
# instead of minting brand-new styles), then overwrite with real values.
$ws.Cells.Item($prevRow, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)

$ws.Cells.Item($prevRow, 5).Copy()
$ws.Cells.Item($newRow, 5).PasteSpecial(-4122)

$ws.Cells.Item($newRow, 1).Value = 76
$ws.Cells.Item($newRow, 2).Value = "armenia"
$ws.Cells.Item($newRow, 3).Value = "premier-league"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45241.5
$ws.Cells.Item($newRow, 6).Value = "Pyunik Yerevan"
$ws.Cells.Item($newRow, 7).Value = 3
$ws.Cells.Item($newRow, 8).Value = "Noah"
$ws.Cells.Item($newRow, 9).Value = 1
$ws.Cells.Item($newRow, 10).Value = 1.83
$ws.Cells.Item($newRow, 11).Value = "10/11/2023 00:12"
$ws.Cells.Item($newRow, 12).Value = 1.56
$ws.Cells.Item($newRow, 13).Value = "11/11/2023 10:23"
$ws.Cells.Item($newRow, 14).Value = 4.35
$ws.Cells.Item($newRow, 15).Value = "10/11/2023 00:12"
$ws.Cells.Item($newRow, 16).Value = 4.85
$ws.Cells.Item($newRow, 17).Value = "11/11/2023 11:04"
$ws.Cells.Item($newRow, 18).Value = 3.18
$ws.Cells.Item($newRow, 19).Value = "10/11/2023 00:12"
$ws.Cells.Item($newRow, 20).Value = 4.78
$ws.Cells.Item($newRow, 21).Value = "11/11/2023 11:04"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/armenia/premier-league/pyunik-yerevan-noah/f1rEyVYN/"
